$wb = $excel.ActiveWorkbook

$wsExp = $wb.Worksheets.Item("exp")
$wsMech = $wb.Worksheets.Item("mech")

# "exp" sheet: new experiment data file (phi05 -> phi1)
$wsExp.Range("A2").Value = "dames_2016_c3h8_rcm_idt_phi1.xlsx"

# "mech" sheet: new mechanism file + name (reduced -> reduced 5)
$wsMech.Range("A2").Value = "nuig1.2_reduced5.cti"
$wsMech.Range("C2").Value = "NUIG1.2, reduced 5"

# widen column C on "mech" to fit the longer mechanism-name text
$wsMech.Columns.Item(3).ColumnWidth = 16.42

# update the cell selections left on each sheet
$wsExp.Range("B12").Select()
$wsMech.Range("B13").Select()
